$d = $word.ActiveDocument

# The document ends with:
#   ... "LOQ4095: Química Geral Experimental (Requisito)"
#   <empty paragraph>
#   <empty paragraph, page-break-before>
#   "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. ..."
#   <empty paragraph>
#   <empty paragraph, page-break-before>
#
# The edit removes the blank paragraph, the page-break paragraph, and the
# copyright paragraph that sit right after the "LOQ4095" requirement line,
# while leaving the final two (identical-looking) blank paragraphs in place.

$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*LOQ4095*") {
        $anchor = $p
    }
}

if ($anchor -eq $null) {
    throw "Could not find the 'LOQ4095' anchor paragraph"
}

$firstToRemove = $anchor.Next()
$lastToRemove = $firstToRemove.Next().Next()

$deleteRange = $d.Range($firstToRemove.Range.Start, $lastToRemove.Range.End)
$deleteRange.Delete()
